$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(1.0143007077573436, 0.6385915078721041, 1.793698772100036, 3.574898503866785)"
$ws.Range("C2").Value = "NIG(1.6688922422492218, 1.321989838963093, 3.5409496724014566, 5.600176281744595)"
$ws.Range("D2").Value = "NIG(0.6989422650214836, 0.5132186415647302, 0.9363410084649437, 2.9267378691672605)"
$ws.Range("E2").Value = "JSU(-1.601379415942842, 1.1742660808338314, 1.4335939892558642, 4.407526933365681)"
